$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.185.96"
$ws.Range("E2").Value = "  -0.09%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.852.68"
$ws.Range("E3").Value = "  -0.57%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.16%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "235.54"
$ws.Range("E5").Value = "  +0.24%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  +0.22%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4712"
$ws.Range("E7").Value = "  +0.99%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2895"
$ws.Range("E8").Value = "  +2.30%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06554"
$ws.Range("E9").Value = "  +0.29%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.71"
$ws.Range("E10").Value = "  +1.27%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07949"
$ws.Range("E11").Value = "  +1.20%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "97.46"
$ws.Range("E12").Value = "  +0.02%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.864.80"
$ws.Range("E13").Value = "  +0.04%  "

$ws.Range("E14").Value = "  -0.34%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6766"
$ws.Range("E15").Value = "  +0.58%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "267.74"
$ws.Range("E16").Value = "  -4.51%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.159.17"
$ws.Range("E17").Value = "  -0.16%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.66"
$ws.Range("E18").Value = "  +8.05%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007571"
$ws.Range("E19").Value = "  +4.00%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.002"
$ws.Range("E20").Value = "  +0.16%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.102.02"
$ws.Range("E21").Value = "  -0.48%  "

$ws.Range("E22").Value = "  +0.18%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.240"
$ws.Range("E23").Value = "  -5.05%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.123"
$ws.Range("E24").Value = "  -0.40%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "166.33"
$ws.Range("E25").Value = "  +1.01%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.142"
$ws.Range("E26").Value = "  -0.61%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.81"
$ws.Range("E27").Value = "  -1.80%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.934"
$ws.Range("E28").Value = "  +0.51%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.399"
$ws.Range("E29").Value = "  +1.53%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09875"
$ws.Range("E30").Value = "  +1.85%  "

$ws.Range("E31").Value = "  -0.58%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.291"
$ws.Range("E32").Value = "  -2.89%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.000"
$ws.Range("E33").Value = "  -2.07%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04680"
$ws.Range("E34").Value = "  -0.28%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.122"
$ws.Range("E35").Value = "  +0.71%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6968"
$ws.Range("E36").Value = "  -1.27%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.712"
$ws.Range("E37").Value = "  -0.60%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01869"
$ws.Range("E38").Value = "  +0.82%  "

$ws.Range("E39").Value = "  +3.09%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.322"
$ws.Range("E40").Value = "  +1.56%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "73.39"
$ws.Range("E41").Value = "  +0.17%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.916"
$ws.Range("E42").Value = "  -1.05%  "

$ws.Range("E43").Value = "  +0.23%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8367"
$ws.Range("E44").Value = "  -1.31%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "103.28"
$ws.Range("E45").Value = "  -0.66%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4135"
$ws.Range("E46").Value = "  -0.72%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "943.77"
$ws.Range("E47").Value = "  +1.04%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.083"
$ws.Range("E48").Value = "  -0.75%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.970"
$ws.Range("E49").Value = "  -3.20%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "33.90"
$ws.Range("E50").Value = "  -0.56%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05652"
$ws.Range("E51").Value = "  +0.45%  "
